$wb = $excel.ActiveWorkbook

# Sheet ALC, row 33 (G33=5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 492.8
$ws.Range("I33").Value = 478.14285
$ws.Range("K33").Value = 478.14285
$ws.Range("M33").Value = -249.14285

# Sheet ALC, row 113 (G113=27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 28574342
$ws.Range("I113").Value = 50002476
$ws.Range("J113").Value = 3499.6667
$ws.Range("K113").Value = 50002476
$ws.Range("L113").Value = 3499.6667
$ws.Range("M113").Value = -49999222
$ws.Range("N113").Value = -10007.6667

# Sheet ALC, row 137 (G137=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1638.6666
$ws.Range("I137").Value = 1649.3334
$ws.Range("J137").Value = 1630.1333
$ws.Range("K137").Value = 4948.0002
$ws.Range("L137").Value = 4890.3999
$ws.Range("M137").Value = -2398.0002
$ws.Range("N137").Value = -9990.3999

# Sheet ALC, row 138 (G138=44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3348.2898
$ws.Range("J138").Value = 3248.2388
$ws.Range("L138").Value = 9744.716400000001
$ws.Range("N138").Value = -20024.7164

# Sheet ARM, row 32 (G32=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22754.254
$ws.Range("I32").Value = 14042.462
$ws.Range("J32").Value = 43989.25
$ws.Range("K32").Value = 14042.462
$ws.Range("L32").Value = 43989.25
$ws.Range("M32").Value = -13755.462
$ws.Range("N32").Value = -44563.25

# Sheet ARM, row 80 (G80=10667)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 37400
$ws.Range("J80").Value = 37400
$ws.Range("L80").Value = 37400
$ws.Range("N80").Value = -39396

# Sheet ARM, row 83 (G83=10667)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 37400
$ws.Range("J83").Value = 37400
$ws.Range("L83").Value = 112200
$ws.Range("N83").Value = -122184

# Sheet ARM, row 88 (G88=12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1798.25
$ws.Range("I88").Value = 1379.2
$ws.Range("J88").Value = 2496.6667
$ws.Range("K88").Value = 1379.2
$ws.Range("L88").Value = 2496.6667
$ws.Range("M88").Value = -973.2
$ws.Range("N88").Value = -3308.6667

# Sheet ARM, row 91 (G91=12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1798.25
$ws.Range("I91").Value = 1379.2
$ws.Range("J91").Value = 2496.6667
$ws.Range("K91").Value = 1379.2
$ws.Range("L91").Value = 2496.6667
$ws.Range("M91").Value = 24.79999999999995
$ws.Range("N91").Value = -5304.6667

# Sheet ARM, row 132 (G132=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2474.5334
$ws.Range("I132").Value = 2170.8635
$ws.Range("K132").Value = 6512.5905
$ws.Range("M132").Value = -3982.5905

# Sheet BSM, row 24 (G24=2420)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 498
$ws.Range("I24").Value = 498
$ws.Range("K24").Value = 498
$ws.Range("M24").Value = -263

# Sheet BSM, row 29 (G29=2318)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 14011.667
$ws.Range("I29").Value = 1999
$ws.Range("K29").Value = 1999
$ws.Range("M29").Value = -1710

# Sheet BSM, row 58 (G58=43234)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 33689
$ws.Range("J58").Value = 39611.25
$ws.Range("L58").Value = 39611.25
$ws.Range("N58").Value = -40199.25

# Sheet BSM, row 86 (G86=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4511.05
$ws.Range("I86").Value = 4365.5293
$ws.Range("J86").Value = 5335.6665
$ws.Range("K86").Value = 4365.5293
$ws.Range("L86").Value = 5335.6665
$ws.Range("M86").Value = -3242.5293
$ws.Range("N86").Value = -7581.6665

# Sheet BSM, row 89 (G89=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4511.05
$ws.Range("I89").Value = 4365.5293
$ws.Range("J89").Value = 5335.6665
$ws.Range("K89").Value = 21827.6465
$ws.Range("L89").Value = 26678.3325
$ws.Range("M89").Value = -16211.6465
$ws.Range("N89").Value = -37910.3325

# Sheet BSM, row 134 (G134=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3981.1943
$ws.Range("I134").Value = 877.625
$ws.Range("J134").Value = 6464.05
$ws.Range("K134").Value = 2632.875
$ws.Range("L134").Value = 19392.15
$ws.Range("M134").Value = -97.875
$ws.Range("N134").Value = -24462.15

# Sheet CRP, row 31 (G31=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1398.8125
$ws.Range("I31").Value = 1398.8125
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1398.8125
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1103.8125
$ws.Range("N31").Value = $null

# Sheet CRP, row 34 (G34=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1398.8125
$ws.Range("I34").Value = 1398.8125
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1398.8125
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1196.8125
$ws.Range("N34").Value = $null

# Sheet CUL, row 11 (G11=4745)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 163381.11
$ws.Range("I11").Value = 193829.44
$ws.Range("J11").Value = 990
$ws.Range("K11").Value = 581488.3200000001
$ws.Range("L11").Value = 2970
$ws.Range("M11").Value = -581348.3200000001
$ws.Range("N11").Value = -3250

# Sheet CUL, row 12 (G12=4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 211.53847
$ws.Range("I12").Value = 349.66666
$ws.Range("J12").Value = 170.1
$ws.Range("K12").Value = 1048.99998
$ws.Range("L12").Value = 510.3
$ws.Range("M12").Value = -875.9999800000001
$ws.Range("N12").Value = -856.3

# Sheet CUL, row 31 (G31=4710)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2101
$ws.Range("J31").Value = 2101
$ws.Range("L31").Value = 6303
$ws.Range("N31").Value = -6879

# Sheet CUL, row 113 (G113=27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 684.0833
$ws.Range("I113").Value = 550
$ws.Range("K113").Value = 1650
$ws.Range("M113").Value = 520

# Sheet CUL, row 117 (G117=27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 759
$ws.Range("I117").Value = 624.3
$ws.Range("K117").Value = 1872.9
$ws.Range("M117").Value = 1569.1

# Sheet CUL, row 121 (G121=27878)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 938.75
$ws.Range("J121").Value = 1177.5
$ws.Range("L121").Value = 3532.5
$ws.Range("N121").Value = -6152.5

# Sheet CUL, row 131 (G131=36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 18210052
$ws.Range("I131").Value = 100000430
$ws.Range("J131").Value = 34411.2
$ws.Range("K131").Value = 300001290
$ws.Range("L131").Value = 103233.6
$ws.Range("M131").Value = -299996250
$ws.Range("N131").Value = -113313.6

# Sheet CUL, row 137 (G137=44088)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 34098660
$ws.Range("I137").Value = 68184080
$ws.Range("J137").Value = 13242.363
$ws.Range("K137").Value = 204552240
$ws.Range("L137").Value = 39727.089
$ws.Range("M137").Value = -204547140
$ws.Range("N137").Value = -49927.089

# Sheet GSM, row 25 (G25=2687)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = $null

# Sheet GSM, row 132 (G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6962.231
$ws.Range("I132").Value = 7774.6
$ws.Range("K132").Value = 23323.8
$ws.Range("M132").Value = -20793.8

# Sheet LTW, row 7 (G7=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2753.7693
$ws.Range("I7").Value = 2528.5715
$ws.Range("K7").Value = 2528.5715
$ws.Range("M7").Value = -2416.5715

# Sheet LTW, row 40 (G40=36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3088.6667
$ws.Range("I40").Value = 3114
$ws.Range("K40").Value = 3114
$ws.Range("M40").Value = -2978

# Sheet LTW, row 126 (G126=36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2753.7693
$ws.Range("I126").Value = 2528.5715
$ws.Range("K126").Value = 7585.7145
$ws.Range("M126").Value = -5115.7145

# Sheet WVR, row 62 (G62=12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100003200
$ws.Range("I62").Value = 125002750
$ws.Range("K62").Value = 125002750
$ws.Range("M62").Value = -125002126

# Sheet WVR, row 65 (G65=12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 100003200
$ws.Range("I65").Value = 125002750
$ws.Range("K65").Value = 625013750
$ws.Range("M65").Value = -625010630

# Sheet WVR, row 122 (G122=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8066207.5
$ws.Range("J122").Value = 799.5
$ws.Range("L122").Value = 2398.5
$ws.Range("N122").Value = -7298.5
